$wb = $excel.ActiveWorkbook

# --- Create the new "Logging" sheet as a copy of "Localizer", placed right after it ---
$loc = $wb.Worksheets.Item("Localizer")
$loc.Copy([System.Type]::Missing, $loc)
$log = $wb.Worksheets.Item($loc.Index + 1)
$log.Name = "Logging"

# --- Update the textual content of the new sheet ---
$log.Range("C2").Value = "Logging"

$log.Range("C4").Value = "Can log information"
$log.Range("C5").Value = "Can log warning"
$log.Range("C6").Value = "Can log error"
$log.Range("C7").Value = "Can log fatal error"

$log.Range("G4").Value = "Testing.Logging.LogsHandlerTests.LogInformationTest()"
$log.Range("G5").Value = "Testing.Logging.LogsHandlerTests.LogWarningTest()"
$log.Range("G6").Value = "Testing.Logging.LogsHandlerTests.LogErrorTest()"
$log.Range("G7").Value = "Testing.Logging.LogsHandlerTests.LogFatalErrorTest()"

# --- Fix up formatting differences vs. the copied "Localizer" layout ---
# G3 on "Logging" uses the plain header border style (same as H3), not the
# "Localizer"-specific one.
$log.Range("H3").Copy()
$log.Range("G3").PasteSpecial(-4122) | Out-Null
$log.Application.CutCopyMode = $false

# Row 5 keeps the default row height on "Logging" (no wrapped/tall row).
$log.Rows.Item(5).AutoFit() | Out-Null

# Column widths differ slightly from "Localizer" on the new sheet.
$log.Columns.Item(3).ColumnWidth = 27.75
$log.Columns.Item(7).ColumnWidth = 49.75

# --- Selections: "Logging" becomes the active sheet/tab with G10 selected ---
$loc.Activate()
$loc.Range("B2:H7").Select()

$log.Activate()
$log.Range("G10").Select()
